$d = $word.ActiveDocument

$wNs  = 'xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"'
$w14Ns = 'xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml"'

# ---------------------------------------------------------------------------
# Paragraph "Crystal": remove the spellStart/spellEnd proofErr markers that
# surround the run. The run content itself (text + formatting) is unchanged.
# ---------------------------------------------------------------------------
$rCrystal = $d.Content
$foundCrystal = $rCrystal.Find.Execute("Crystal", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if ($foundCrystal -and $rCrystal.Text -eq "Crystal") {
    $xmlCrystal = '<w:p ' + $wNs + ' ' + $w14Ns + ' w14:paraId="7ED6FF47" w14:textId="77777777" w:rsidR="00ED4DB3" w:rsidRPr="00ED4DB3" w:rsidRDefault="00ED4DB3" w:rsidP="00ED4DB3">' +
                  '<w:pPr><w:pStyle w:val="Prrafodelista"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="4"/></w:numPr><w:tabs><w:tab w:val="left" w:pos="8730"/></w:tabs><w:rPr><w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/><w:lang w:val="es-ES"/></w:rPr></w:pPr>' +
                  '<w:r w:rsidRPr="00ED4DB3"><w:rPr><w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/></w:rPr><w:t>Crystal</w:t></w:r>' +
                  '</w:p>'
    $rCrystal.InsertXML($xmlCrystal)
}

# ---------------------------------------------------------------------------
# Paragraph "Lean " + "Development": merge the two runs into a single run
# with the text "Lean Development" and drop the proofErr markers that used
# to wrap the second ("Development") run.
# ---------------------------------------------------------------------------
$rLean = $d.Content
$foundLean = $rLean.Find.Execute("Lean Development", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if ($foundLean -and $rLean.Text -eq "Lean Development") {
    $xmlLean = '<w:p ' + $wNs + ' ' + $w14Ns + ' w14:paraId="190C4F17" w14:textId="77777777" w:rsidR="00ED4DB3" w:rsidRPr="00ED4DB3" w:rsidRDefault="00ED4DB3" w:rsidP="00ED4DB3">' +
               '<w:pPr><w:pStyle w:val="Prrafodelista"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="4"/></w:numPr><w:tabs><w:tab w:val="left" w:pos="8730"/></w:tabs><w:rPr><w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/><w:lang w:val="es-ES"/></w:rPr></w:pPr>' +
               '<w:r w:rsidRPr="00ED4DB3"><w:rPr><w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/></w:rPr><w:t>Lean Development</w:t></w:r>' +
               '</w:p>'
    $rLean.InsertXML($xmlLean)
}

Write-Output "done"
